$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Title

$titleShape.TextFrame.TextRange.Text = "Implementacja Eliminacji Gaussa-Jordana w różnych wariantach"
$titleShape.Height = 245.208865
